$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bisseccao: epsilon tightened from 1e-4 to 1e-6, results refined
$ws.Range("B2").Value = "a=1.0, b=2.0, ε=1e-06"
$ws.Range("C2").Value = 1.324718
$ws.Range("D2").Value = [double]"1.75958307e-07"
$ws.Range("E2").Value = [double]"9.53674316e-07"
$ws.Range("F2").Value = 20

# Row 3 - Posicao Falsa: epsilon tightened from 1e-4 to 1e-6, results refined
$ws.Range("B3").Value = "a=1.0, b=2.0, ε=1e-06"
$ws.Range("C3").Value = 1.32471787
$ws.Range("D3").Value = [double]"-3.51552278e-07"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 18

# Row 4 - Ponto Fixo: epsilon tightened from 1e-4 to 1e-6, results refined
$ws.Range("B4").Value = "x0=1.2, ε=1e-06"
$ws.Range("C4").Value = 1.32471774
$ws.Range("D4").Value = [double]"-9.21852971e-07"
$ws.Range("E4").Value = [double]"9.21852971e-07"
$ws.Range("F4").Value = 8

# Row 5 - Newton: starting point changed to 0.5, epsilon tightened to 1e-6
$ws.Range("B5").Value = "x0=0.5, ε=1e-06"
$ws.Range("C5").Value = 1.324718
$ws.Range("D5").Value = [double]"1.74374144e-07"
$ws.Range("E5").Value = 0.000209457449
$ws.Range("F5").Value = 19

# Row 6 - Secante: starting points changed to 0.0/0.5, epsilon tightened to 1e-6
# C6 becomes a genuine number; D6/E6 stay textual, so force Text format first
$ws.Range("B6").Value = "x0=0.0, x1=0.5, ε=1e-06"
$ws.Range("C6").Value = 1.32471795
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "-4.34057552e-08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.19166137e-05"
$ws.Range("F6").Value = 26
